# Minor wording improvements in the template
#
# This script rewrites the instructional text cells on Sheet1 so that:
#   - Step 1 now explains adding "x" in the "My name (x)" column (new wording)
#   - Step 2 now references rating in the "Rating" column (new wording)
#   - Step 4 now refers to the "Comments" column (quoted, new wording)
# and moves the selection highlight to A11 (matching the author's last
# selected cell when saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1 (A2): reworded instruction -------------------------------------
$ws.Range("A2").Value = "Step 1: In the table below, on the row with your name, add ""x"" (without quotes) in column ""My name (x)""."

# --- Step 2 (A4): reworded instruction, keep rich-text emphasis on "your" --
$step2Text = "Step 2: Rate your contribution in the column ""Rating""to the project on the range 1" + [char]0x2013 + "5 according to the following guide:"
$ws.Range("A4").Value = $step2Text

# "your" (characters 14-17) stays bold + underlined, in "Calibri (Body)"
$step2Your = $ws.Range("A4").Characters(14, 4)
$step2Your.Font.Bold = $true
$step2Your.Font.Underline = $true
$step2Your.Font.Name = "Calibri (Body)"

# remainder of the sentence stays bold, regular (non-underlined) "Calibri"
$step2RestLength = $step2Text.Length - 17
$step2Rest = $ws.Range("A4").Characters(18, $step2RestLength)
$step2Rest.Font.Bold = $true

# --- Step 4 (A13): reworded instruction, "Comments" now quoted -------------
$ws.Range("A13").Value = "Step 4: If you have any comments about any person's contribution, write in the ""Comments"" column"

# --- Restore the last active selection (A11) --------------------------------
$ws.Range("A11").Select()
